# "8. Primary Key.pptx" update
#
# 1) Slide 1 ("Rectangle 5" body text box) - the three separate runs describing
#    what a primary key is are merged into a single run (uniformly formatted
#    in Palatino Linotype) and the box is made taller to fit the combined
#    paragraph.
# 2) Slide 5 (the "INSERT" slide) is removed from the deck entirely.

$p = $ppt.ActivePresentation

# --- 1) Slide 1: merge the "Choosing a primary key ..." paragraph ------------
$slide1 = $p.Slides.Item(1)
$bodyShape = $slide1.Shapes.Item("Rectangle 5")
$tr = $bodyShape.TextFrame.TextRange

$fullText = "Choosing a primary key is one of the most important steps in good database design. A primary key is a column that serves a special purpose. A primary key is a special column (or set of combined columns) in a relational database table, that is used to uniquely identify each record. Each database table needs a primary key."

# Write a throwaway value first so the subsequent assignment can't be folded
# back onto the pre-existing (differently formatted) runs - this guarantees
# the paragraph ends up as one single run once we restyle it below.
$tr.Text = "-"
$tr.Text = $fullText
$tr.Font.Name = "Palatino Linotype"

# The box uses "resize shape to fit text" - the extra sentence needs more
# vertical room, so grow the box to match the new wrapped height.
$bodyShape.Height = 94.5141

# --- 2) Drop the 5th slide ("INSERT ... statement") --------------------------
$p.Slides.Item(5).Delete()
